$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Totals")

$ws.Range("B19").Value = 2809
$ws.Range("C19").Value = 66573
$ws.Range("E19").Value = 39180
$ws.Range("F19").Value = 108563

$ws.Range("B20").Select()
